$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the precision of the existing A83 timestamp value.
$ws.Range("A83").Value = 44396.77000582755

# Append new row 84 with the latest retrieved data point.
$ws.Range("A84").Value = 44397.76866720003
$ws.Range("B84").Value = 79821
$ws.Range("C84").Value = 67296
$ws.Range("D84").Value = 3728
$ws.Range("E84").Value = 2195
$ws.Range("F84").Value = 1588
$ws.Range("G84").Value = 20842
$ws.Range("H84").Value = 1632
$ws.Range("I84").Value = 886
$ws.Range("J84").Value = 203

# New date cell should carry the same number format as the rest of column A.
$ws.Range("A84").NumberFormat = $ws.Range("A83").NumberFormat
